$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.721116280972109, 1.5114120164233908, -0.966702854389722, 2.533985683069095)"
$ws.Range("C2").Value = "NIG(2.031888860545584, 1.6774120902836565, 3.4159758905180535, 6.23074190749257)"
$ws.Range("D2").Value = "NIG(0.9163159414801771, 0.6201668573895849, 1.2056742748848843, 3.2365445692898343)"
$ws.Range("E2").Value = "NIG(1.3323388630239843, 1.0274026294851104, 2.740350203058832, 6.479167014603409)"
